$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.323.07"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.933.12"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7497"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "249.06"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.02"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3216"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07119"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7894"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08006"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.937.72"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.380"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.49"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.326.24"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "253.55"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008039"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.805"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.193.07"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.830"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.96%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.40"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.46%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.320"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.55%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.11"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1341"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.356"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.532"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.428"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.147"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05122"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.288"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7494"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.773"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01966"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.807"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.05"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.413"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4517"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.990"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8418"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.19"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.821"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.550"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "984.96"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +11.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.46"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1197"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.45%  "
